$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("October")

# Fill in October's net borrower/lender statistics (previously blank placeholder rows).
# For each library row: B = borrowed, C = lent, D = B - C (net),
# E/F = explanatory text depending on the sign of D, G = ratio "x.xx : 1".

$ws.Range("B2").Value = 1716
$ws.Range("C2").Value = 1294
$ws.Range("D2").Value = 422
$ws.Range("E2").Value = "We borrowerd more than we lent"
$ws.Range("G2").Value = "1.33 : 1"

$ws.Range("B3").Value = 712
$ws.Range("C3").Value = 416
$ws.Range("D3").Value = 296
$ws.Range("E3").Value = "We borrowerd more than we lent"
$ws.Range("G3").Value = "1.71 : 1"

$ws.Range("B4").Value = 1213
$ws.Range("C4").Value = 1459
$ws.Range("D4").Value = -246
$ws.Range("F4").Value = "We lent more than we borrowed"
$ws.Range("G4").Value = "0.83 : 1"

$ws.Range("B5").Value = 33
$ws.Range("C5").Value = 157
$ws.Range("D5").Value = -124
$ws.Range("F5").Value = "We lent more than we borrowed"
$ws.Range("G5").Value = "0.21 : 1"

$ws.Range("B6").Value = 1258
$ws.Range("C6").Value = 1549
$ws.Range("D6").Value = -291
$ws.Range("F6").Value = "We lent more than we borrowed"
$ws.Range("G6").Value = "0.81 : 1"

$ws.Range("B7").Value = 175
$ws.Range("C7").Value = 191
$ws.Range("D7").Value = -16
$ws.Range("F7").Value = "We lent more than we borrowed"
$ws.Range("G7").Value = "0.92 : 1"

$ws.Range("B8").Value = 181
$ws.Range("C8").Value = 207
$ws.Range("D8").Value = -26
$ws.Range("F8").Value = "We lent more than we borrowed"
$ws.Range("G8").Value = "0.87 : 1"

$ws.Range("B9").Value = 55
$ws.Range("C9").Value = 77
$ws.Range("D9").Value = -22
$ws.Range("F9").Value = "We lent more than we borrowed"
$ws.Range("G9").Value = "0.71 : 1"

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 42
$ws.Range("D10").Value = -42
$ws.Range("F10").Value = "We lent more than we borrowed"
$ws.Range("G10").Value = "0.00 : 1"

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 17
$ws.Range("C12").Value = 12
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = "We borrowerd more than we lent"
$ws.Range("G12").Value = "1.42 : 1"

$ws.Range("B13").Value = 111
$ws.Range("C13").Value = 75
$ws.Range("D13").Value = 36
$ws.Range("E13").Value = "We borrowerd more than we lent"
$ws.Range("G13").Value = "1.48 : 1"

$ws.Range("B14").Value = 192
$ws.Range("C14").Value = 235
$ws.Range("D14").Value = -43
$ws.Range("F14").Value = "We lent more than we borrowed"
$ws.Range("G14").Value = "0.82 : 1"

$ws.Range("B15").Value = 55
$ws.Range("C15").Value = 171
$ws.Range("D15").Value = -116
$ws.Range("F15").Value = "We lent more than we borrowed"
$ws.Range("G15").Value = "0.32 : 1"

$ws.Range("B16").Value = 38
$ws.Range("C16").Value = 134
$ws.Range("D16").Value = -96
$ws.Range("F16").Value = "We lent more than we borrowed"
$ws.Range("G16").Value = "0.28 : 1"

$ws.Range("B17").Value = 701
$ws.Range("C17").Value = 434
$ws.Range("D17").Value = 267
$ws.Range("E17").Value = "We borrowerd more than we lent"
$ws.Range("G17").Value = "1.62 : 1"

$ws.Range("B18").Value = 93
$ws.Range("C18").Value = 160
$ws.Range("D18").Value = -67
$ws.Range("F18").Value = "We lent more than we borrowed"
$ws.Range("G18").Value = "0.58 : 1"

$ws.Range("B19").Value = 615
$ws.Range("C19").Value = 455
$ws.Range("D19").Value = 160
$ws.Range("E19").Value = "We borrowerd more than we lent"
$ws.Range("G19").Value = "1.35 : 1"

$ws.Range("B20").Value = 49
$ws.Range("C20").Value = 69
$ws.Range("D20").Value = -20
$ws.Range("F20").Value = "We lent more than we borrowed"
$ws.Range("G20").Value = "0.71 : 1"

$ws.Range("B21").Value = 557
$ws.Range("C21").Value = 361
$ws.Range("D21").Value = 196
$ws.Range("E21").Value = "We borrowerd more than we lent"
$ws.Range("G21").Value = "1.54 : 1"

$ws.Range("B22").Value = 53
$ws.Range("C22").Value = 102
$ws.Range("D22").Value = -49
$ws.Range("F22").Value = "We lent more than we borrowed"
$ws.Range("G22").Value = "0.52 : 1"

$ws.Range("B23").Value = 635
$ws.Range("C23").Value = 344
$ws.Range("D23").Value = 291
$ws.Range("E23").Value = "We borrowerd more than we lent"
$ws.Range("G23").Value = "1.85 : 1"

$ws.Range("B24").Value = 1623
$ws.Range("C24").Value = 1202
$ws.Range("D24").Value = 421
$ws.Range("E24").Value = "We borrowerd more than we lent"
$ws.Range("G24").Value = "1.35 : 1"

$ws.Range("B25").Value = 169
$ws.Range("C25").Value = 534
$ws.Range("D25").Value = -365
$ws.Range("F25").Value = "We lent more than we borrowed"
$ws.Range("G25").Value = "0.32 : 1"

$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0

$ws.Range("B27").Value = 252
$ws.Range("C27").Value = 221
$ws.Range("D27").Value = 31
$ws.Range("E27").Value = "We borrowerd more than we lent"
$ws.Range("G27").Value = "1.14 : 1"

$ws.Range("B28").Value = 151
$ws.Range("C28").Value = 90
$ws.Range("D28").Value = 61
$ws.Range("E28").Value = "We borrowerd more than we lent"
$ws.Range("G28").Value = "1.68 : 1"

$ws.Range("B29").Value = 701
$ws.Range("C29").Value = 471
$ws.Range("D29").Value = 230
$ws.Range("E29").Value = "We borrowerd more than we lent"
$ws.Range("G29").Value = "1.49 : 1"

$ws.Range("B30").Value = 55
$ws.Range("C30").Value = 33
$ws.Range("D30").Value = 22
$ws.Range("E30").Value = "We borrowerd more than we lent"
$ws.Range("G30").Value = "1.67 : 1"

$ws.Range("B31").Value = 81
$ws.Range("C31").Value = 286
$ws.Range("D31").Value = -205
$ws.Range("F31").Value = "We lent more than we borrowed"
$ws.Range("G31").Value = "0.28 : 1"

$ws.Range("B32").Value = 443
$ws.Range("C32").Value = 610
$ws.Range("D32").Value = -167
$ws.Range("F32").Value = "We lent more than we borrowed"
$ws.Range("G32").Value = "0.73 : 1"

$ws.Range("B33").Value = 367
$ws.Range("C33").Value = 512
$ws.Range("D33").Value = -145
$ws.Range("F33").Value = "We lent more than we borrowed"
$ws.Range("G33").Value = "0.72 : 1"

$ws.Range("B34").Value = 185
$ws.Range("C34").Value = 137
$ws.Range("D34").Value = 48
$ws.Range("E34").Value = "We borrowerd more than we lent"
$ws.Range("G34").Value = "1.35 : 1"

$ws.Range("B35").Value = 910
$ws.Range("C35").Value = 1300
$ws.Range("D35").Value = -390
$ws.Range("F35").Value = "We lent more than we borrowed"
$ws.Range("G35").Value = "0.70 : 1"

$ws.Range("B36").Value = 186
$ws.Range("C36").Value = 572
$ws.Range("D36").Value = -386
$ws.Range("F36").Value = "We lent more than we borrowed"
$ws.Range("G36").Value = "0.33 : 1"

$ws.Range("B37").Value = 615
$ws.Range("C37").Value = 489
$ws.Range("D37").Value = 126
$ws.Range("E37").Value = "We borrowerd more than we lent"
$ws.Range("G37").Value = "1.26 : 1"

$ws.Range("B38").Value = 46
$ws.Range("C38").Value = 186
$ws.Range("D38").Value = -140
$ws.Range("F38").Value = "We lent more than we borrowed"
$ws.Range("G38").Value = "0.25 : 1"

$ws.Range("B39").Value = 35
$ws.Range("C39").Value = 109
$ws.Range("D39").Value = -74
$ws.Range("F39").Value = "We lent more than we borrowed"
$ws.Range("G39").Value = "0.32 : 1"

$ws.Range("B40").Value = 59
$ws.Range("C40").Value = 147
$ws.Range("D40").Value = -88
$ws.Range("F40").Value = "We lent more than we borrowed"
$ws.Range("G40").Value = "0.40 : 1"

$ws.Range("B41").Value = 20
$ws.Range("C41").Value = 41
$ws.Range("D41").Value = -21
$ws.Range("F41").Value = "We lent more than we borrowed"
$ws.Range("G41").Value = "0.49 : 1"

$ws.Range("B42").Value = 8
$ws.Range("C42").Value = 56
$ws.Range("D42").Value = -48
$ws.Range("F42").Value = "We lent more than we borrowed"
$ws.Range("G42").Value = "0.14 : 1"

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0

$ws.Range("B44").Value = 98
$ws.Range("C44").Value = 103
$ws.Range("D44").Value = -5
$ws.Range("F44").Value = "We lent more than we borrowed"
$ws.Range("G44").Value = "0.95 : 1"

$ws.Range("B45").Value = 43
$ws.Range("C45").Value = 194
$ws.Range("D45").Value = -151
$ws.Range("F45").Value = "We lent more than we borrowed"
$ws.Range("G45").Value = "0.22 : 1"

$ws.Range("B46").Value = 626
$ws.Range("C46").Value = 624
$ws.Range("D46").Value = 2
$ws.Range("E46").Value = "We borrowerd more than we lent"
$ws.Range("G46").Value = "1.00 : 1"

$ws.Range("B47").Value = 1154
$ws.Range("C47").Value = 635
$ws.Range("D47").Value = 519
$ws.Range("E47").Value = "We borrowerd more than we lent"
$ws.Range("G47").Value = "1.82 : 1"

$ws.Range("B48").Value = 242
$ws.Range("C48").Value = 618
$ws.Range("D48").Value = -376
$ws.Range("F48").Value = "We lent more than we borrowed"
$ws.Range("G48").Value = "0.39 : 1"

$ws.Range("B49").Value = 621
$ws.Range("C49").Value = 270
$ws.Range("D49").Value = 351
$ws.Range("E49").Value = "We borrowerd more than we lent"
$ws.Range("G49").Value = "2.30 : 1"

$ws.Range("B50").Value = 807
$ws.Range("C50").Value = 565
$ws.Range("D50").Value = 242
$ws.Range("E50").Value = "We borrowerd more than we lent"
$ws.Range("G50").Value = "1.43 : 1"

$ws.Range("B51").Value = 268
$ws.Range("C51").Value = 167
$ws.Range("D51").Value = 101
$ws.Range("E51").Value = "We borrowerd more than we lent"
$ws.Range("G51").Value = "1.60 : 1"

$ws.Range("B52").Value = 374
$ws.Range("C52").Value = 449
$ws.Range("D52").Value = -75
$ws.Range("F52").Value = "We lent more than we borrowed"
$ws.Range("G52").Value = "0.83 : 1"

$ws.Range("B53").Value = 203
$ws.Range("C53").Value = 241
$ws.Range("D53").Value = -38
$ws.Range("F53").Value = "We lent more than we borrowed"
$ws.Range("G53").Value = "0.84 : 1"

$ws.Range("B54").Value = 51
$ws.Range("C54").Value = 257
$ws.Range("D54").Value = -206
$ws.Range("F54").Value = "We lent more than we borrowed"
$ws.Range("G54").Value = "0.20 : 1"

$ws.Range("B55").Value = 407
$ws.Range("C55").Value = 196
$ws.Range("D55").Value = 211
$ws.Range("E55").Value = "We borrowerd more than we lent"
$ws.Range("G55").Value = "2.08 : 1"
